# "Prix Spot" sheet: a new daily column "04-nov" is inserted right before
# the column that currently holds "01-oct." (column DI), pushing the
# 01-oct..31-oct columns one slot to the right (DI:EM -> DJ:EN).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Prix Spot")

# Insert a blank column before DI; everything from DI onward shifts right.
$ws.Range("DI1").EntireColumn.Insert()

# Header for the freshly inserted column.
$ws.Range("DI1").Value2 = "04-nov"

# The new column has no data yet for this date, shown as "-" like the
# neighbouring not-yet-available "xx-nov" columns (DF:DH).
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 113).Value2 = "-"
}

# "Gaz" sheet: last two rows' prices were revised down from 29.8 to 29.3.
$wsGaz = $wb.Worksheets.Item("Gaz")
$wsGaz.Range("B140").Value2 = 29.3
$wsGaz.Range("B141").Value2 = 29.3
